# Update BGDP (India:US GDP per capita adjustment) on the "About" sheet.
# This single input drives formulas on the "AVMC-passenger" and
# "AVMC-freight" sheets (both reference About!$A$75), so updating it here
# and letting Excel recalculate reproduces all of the downstream value
# changes seen in the diff.

$wb = $excel.ActiveWorkbook

$aboutSheet = $wb.Worksheets.Item("About")
$aboutSheet.Range("A75").Value = 0.03878298458735905

# Recalculate so dependent formulas on AVMC-passenger / AVMC-freight pick up
# the new value.
$excel.Calculate()

# Restore/update the view state on the About sheet to match the author's
# scroll position and selection after editing near the bottom of the sheet.
$aboutSheet.Activate()
$aboutWindow = $excel.ActiveWindow
$aboutWindow.ScrollRow = 64
$aboutSheet.Range("A75").Select()
